$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added "NroPropuesta" column
$ws.Range("S1").Value = "NroPropuesta"

# Duplicate row 2's formatting down into row 3 (new proposal record)
$ws.Range("A2:R2").Copy()
$ws.Range("A3:R3").PasteSpecial(-4122)

# Populate the new row with the proposal's data
$ws.Range("A3").Value = "19499545"
$ws.Range("B3").Value = "CREDITO EMPRESARIAL"
$ws.Range("C3").Value = "CREDITOS PYMES"
$ws.Range("D3").Value = "REFINANCIACION ESPECIAL"
$ws.Range("E3").Value = "NORMAL"
$ws.Range("F3").Value = "SIN PROMOCION"
$ws.Range("G3").Value = "TipoOperacion"
$ws.Range("H3").Value = "20"
$ws.Range("I3").Value = "20"
$ws.Range("J3").Value = "Fija Vencida"
$ws.Range("K3").Value = "Cronograma Pagos"
$ws.Range("L3").Value = "Fecha Fija"
$ws.Range("M3").Value = "5"
$ws.Range("N3").Value = "6"
$ws.Range("O3").Value = "30"
$ws.Range("P3").Value = "EFECTIVO"
$ws.Range("Q3").Value = "prueba de nueva propuesta"
$ws.Range("R3").Value = "prueba de nueva propuesta"

# New formatting-only cell further down the sheet (underlined, default font)
$ws.Range("N9").Value = ""
$ws.Range("N9").Font.Underline = $true

# Column S width to match the new column
$ws.Columns.Item(19).ColumnWidth = 12.5

# Update selection/view to match the authored state
$ws.Range("F17").Select()
